# Working on login screen CSS: add "Front end" build-out notes to the
# logic map (master template / nav / footer / POST request row) and add a
# blank Sheet2 for upcoming front-end work.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Re-colour a few cells that used the old "theme 9" highlight fill to
#     the built-in "Good" (green) cell style instead -------------------
$ws1.Range("G30").Style = "Good"

# --- Rows 36-44: give a handful of previously unstyled label cells the
#     "Good" style, matching the rest of that block -------------------
$ws1.Range("D36").Style = "Good"
$ws1.Range("E37").Style = "Good"
$ws1.Range("F38").Style = "Good"
$ws1.Range("G39").Style = "Good"
$ws1.Range("C44").Style = "Good"

# --- Insert 14 new rows before the old row 52 ("Register new user" block)
#     to make room for a new "Front end" section ------------------------
$ws1.Rows("50:63").Insert()

# New "Front end" section content (shared strings are appended in this
# exact order so they land at indices 73-77, matching Front end, Master
# template, POST request to select reservation, Nav, Footer).
$ws1.Range("A50").Value = "Front end"
$ws1.Range("B51").Value = "Master template"
$ws1.Range("B57").Value = "POST request to select reservation"
$ws1.Range("B52").Value = "Nav"
$ws1.Range("B53").Value = "Footer"

# --- Add a new, empty Sheet2 after Sheet1 ------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"
[void]$ws2.Range("B4").Select()

# --- Restore Sheet1 as the active sheet/selection ----------------------
[void]$ws1.Activate()
[void]$ws1.Range("D51").Select()
$excel.ActiveWindow.ScrollRow = 25
